# Weekly Fruit/Vegetable price update.
# A new daily price record (row 288) is inserted into the "Naranja" sheet,
# pushing all subsequent rows (old 288..367) down by one (new 289..368).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 288; existing rows 288-367 shift down to 289-368.
$ws.Rows.Item(288).EntireRow.Insert()

# Populate the newly inserted row with the new price record.
$ws.Range("A288").Value = 4
$ws.Range("B288").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C288").Value = "Los Lagos"
$ws.Range("D288").Value = 44627
$ws.Range("E288").Value = 10
$ws.Range("F288").Value = "Fruta"
$ws.Range("G288").Value = 100102
$ws.Range("H288").Value = "Cítricos"
$ws.Range("I288").Value = 100102005
$ws.Range("J288").Value = "Naranja"
$ws.Range("K288").Value = "Valencia"
$ws.Range("L288").Value = "Primera"
$ws.Range("M288").Value = 500
$ws.Range("N288").Value = 16000
$ws.Range("O288").Value = 16000
$ws.Range("P288").Value = 16000
$ws.Range("Q288").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R288").Value = "Región de O'Higgins"
$ws.Range("S288").Value = 1067
$ws.Range("T288").Value = 15
